# g1.1 e g1.3 - reestruturação do arquivo
# Reshape the "wide" table (UF/PIB/Rank in A:C, UF/Variação/Rank in D:F)
# into a "long" table (UF/Valor/Rank/Categoria in A:D), stacking the
# "Variação (%) 2022/2010" block underneath the "PIB 2022 Deflacionado"
# block and tagging each with a new "Categoria" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Grab the second block (UF / Variação (%) 2022/2010 / Rank) that
#    currently lives in D2:F10 before anything else gets overwritten.
$variacao = $ws.Range("D2:F10").Value()

# 2. Re-label the header row: B1 "PIB 2022 Deflacionado" -> "Valor",
#    D1 "UF" -> "Categoria"; the old E1/F1 headers go away with the
#    column delete below.
$ws.Range("B1").Value = "Valor"
$ws.Range("D1").Value = "Categoria"

# 3. Tag the existing PIB rows (2-10) with their category in column D.
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 4).Value = "PIB 2022 Deflacionado"
}

# 4. Move the captured Variação block down into A11:C19 ...
$ws.Range("A11:C19").Value = $variacao

# 5. ... and tag those rows with their category in column D.
for ($r = 11; $r -le 19; $r++) {
    $ws.Cells.Item($r, 4).Value = "Variação (%) 2022/2010"
}

# 6. Drop the now-unused E:F columns entirely so the sheet (and its
#    dimension) shrinks back down to A1:D19.
$ws.Range("E1:F19").Delete()
